# Agenda update:
#  - Merge the old "Intro" (8am-8:45am) and "Code Smells" (8:45am-9:45am)
#    sessions into a single "8am-9:45am" session, delivered jointly by
#    Ryan & Cory.
#  - Delete the now-redundant "Code Smells" row (old row 3) so the
#    remaining rows shift up.
#  - Rename "Code Challenges" to "Finding Code Smells" and hand it to
#    Patrick.
#  - Rename "Refactoring practice and assistance" to "Refactoring
#    techniques" and hand it to Patrick.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Merge the 8am slot: new combined time / description.
$ws.Range("A2").Value = "8am-9:45am"
$ws.Range("B2").Value = "Intro & Calories Calculator App"

# Remove the old "Code Smells" row entirely; everything below shifts up.
$ws.Rows.Item(3).Delete()

# "Code Challenges" -> "Finding Code Smells", presented by Patrick.
# (after the deletion this content now lives on row 4)
$ws.Range("B4").Value = "Finding Code Smells"
$ws.Range("C4").Value = "Patrick"

# Refactoring session renamed, now presented by Patrick.
# (after the deletion this content now lives on row 7)
$ws.Range("B7").Value = "Refactoring techniques"
$ws.Range("C7").Value = "Patrick"

# Finish the merged 8am slot with its joint presenter.
$ws.Range("C2").Value = "Ryan & Cory"

# Column widths were manually tweaked to fit the new text.
$ws.Columns.Item(1).ColumnWidth = 16.42578125
$ws.Columns.Item(2).ColumnWidth = 28.5703125

# Selection moved off the table when the edits were done.
$ws.Range("I10").Select()
